{"js": "// Split the long \"Programa\" paragraphs (Portuguese and English) into\n// multiple lines by inserting manual line breaks (<w:br/>) at the\n// sentence boundaries that separate the three topical chunks, matching\n// the author's edit.\n//\n// Strategy: locate a short, unique anchor string that straddles each\n// split point, then replace it with the same text but with a vertical\n// tab character (U+000B) inserted at the split point. Word's Office.js\n// engine renders U+000B inside inserted text as a manual line break\n// (<w:br/>), splitting the run exactly like the target diff.\n\nconst body = context.document.body;\n\nasync function insertLineBreakAt(needle, splitIndex) {\n  const results = body.search(needle, { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length !== 1) {\n    throw new Error(\n      `Expected exactly one match for \"${needle}\", found ${results.items.length}`\n    );\n  }\n\n  const replacement = needle.slice(0, splitIndex) + \"\\u000b\" + needle.slice(splitIndex);\n  results.items[0].insertText(replacement, Word.InsertLocation.replace);\n  await context.sync();\n}\n\n// --- Portuguese paragraph -------------------------------------------------\n// ...fundamentos.| Atom\u00edstica...\nawait insertLineBreakAt(\n  \"cia de materiais; fundamentos.Atom\u00edstica e a vis\u00e3o moderna d\",\n  \"cia de materiais; fundamentos.\".length\n);\n\n// ...liga\u00e7\u00f5es met\u00e1licas.| Materiais policristalinos...\nawait insertLineBreakAt(\n  \" i\u00f4nicas e liga\u00e7\u00f5es met\u00e1licas.Materiais policristalinos e mo\",\n  \" i\u00f4nicas e liga\u00e7\u00f5es met\u00e1licas.\".length\n);\n\n// --- English paragraph -----------------------------------------------------\n// ...fundamentals.| Atomistics...\nawait insertLineBreakAt(\n  \"terials science; fundamentals.Atomistics and the modern view\",\n  \"terials science; fundamentals.\".length\n);\n\n// ...metallic bonds.| Polycrystalline...\nawait insertLineBreakAt(\n  \"onic bonds and metallic bonds.Polycrystalline and monocrysta\",\n  \"onic bonds and metallic bonds.\".length\n);\n", "ps1": "# Split the long \"Programa\" paragraphs (Portuguese and English) into\n# multiple lines by inserting manual line breaks at the sentence\n# boundaries that separate the three topical chunks, matching the\n# author's edit. A manual line break is represented in OOXML as\n# <w:br/>, which Word's Find/Replace produces when the replacement\n# text contains the special \"^l\" sequence.\n\n$d = $word.ActiveDocument\n\nfunction Insert-LineBreak($findText, $replaceText) {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $findText\n    $find.Replacement.Text = $replaceText\n    $find.Forward = $true\n    $find.Wrap = 0            # wdFindStop - do not wrap around / do not loop forever\n    $find.Format = $false\n    $find.MatchCase = $true\n    $find.MatchWholeWord = $false\n    $find.MatchWildcards = $false\n    $find.MatchSoundsLike = $false\n    $find.MatchAllWordForms = $false\n\n    $ok = $find.Execute(\n        [ref]$findText,   # FindText\n        [ref]$true,       # MatchCase\n        [ref]$false,      # MatchWholeWord\n        [ref]$false,      # MatchWildcards\n        [ref]$false,      # MatchSoundsLike\n        [ref]$false,      # MatchAllWordForms\n        [ref]$true,       # Forward\n        [ref]0,           # Wrap (wdFindStop)\n        [ref]$false,      # Format\n        $replaceText,     # ReplaceWith\n        [ref]2            # Replace (wdReplaceOne)\n    )\n\n    if (-not $ok) {\n        throw \"Find/Replace failed for '$findText'\"\n    }\n}\n\n# --- Portuguese paragraph ----------------------------------------------\n# ...fundamentos.| Atom\u00edstica...\nInsert-LineBreak `\n    \"ci\u00eancia de materiais; fundamentos.Atom\u00edstica\" `\n    \"ci\u00eancia de materiais; fundamentos.^lAtom\u00edstica\"\n\n# ...liga\u00e7\u00f5es met\u00e1licas.| Materiais policristalinos...\nInsert-LineBreak `\n    \"liga\u00e7\u00f5es i\u00f4nicas e liga\u00e7\u00f5es met\u00e1licas.Materiais\" `\n    \"liga\u00e7\u00f5es i\u00f4nicas e liga\u00e7\u00f5es met\u00e1licas.^lMateriais\"\n\n# --- English paragraph ---------------------------------------------------\n# ...fundamentals.| Atomistics...\nInsert-LineBreak `\n    \"materials science; fundamentals.Atomistics\" `\n    \"materials science; fundamentals.^lAtomistics\"\n\n# ...metallic bonds.| Polycrystalline...\nInsert-LineBreak `\n    \"ionic bonds and metallic bonds.Polycrystalline\" `\n    \"ionic bonds and metallic bonds.^lPolycrystalline\"\n"}
